$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 13:20"

# Row 12 (Suiza)
$ws.Range("B12").Value = 16176
$ws.Range("C12").Value = 254
$ws.Range("E12").Value = 13980
$ws.Range("G12").Value = 14
$ws.Range("H12").Value = 373

# Row 13 (Belgica)
$ws.Range("D13").Value = 1696
$ws.Range("E13").Value = 10374

# Row 16 (Austria)
$ws.Range("B16").Value = 9876
$ws.Range("C16").Value = 258
$ws.Range("E16").Value = 8653

# Row 27 (Dinamarca)
$ws.Range("E27").Value = 2724
$ws.Range("F27").Value = 145
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 90

# Row 86 (Jordania)
$ws.Range("F86").Value = 5

# Row 92 (Vietnam)
$ws.Range("B92").Value = 207
$ws.Range("C92").Value = 3
$ws.Range("E92").Value = 152

# Row 98 (Malta)
$ws.Range("F98").Value = 2
